$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: the "Förändrad" (changed) date in column C moved from 2023-09-13
# (45182) to 2023-09-15 (45184) for every existing data row (2..348).
$ws.Range("C2:C348").Value = 45184

# Step 2: row 348 picks up an explicit row-height marker once a new row
# gets appended below it.
$ws.Rows.Item(348).RowHeight = 15

# Step 3: append the new record as row 349.
$ws.Range("A349").Value = "A 42976-2023"

$ws.Range("B349").Value = 45182
$ws.Range("B349").NumberFormat = "YYYY-MM-DD"

$ws.Range("C349").Value = 45184
$ws.Range("C349").NumberFormat = "YYYY-MM-DD"

$ws.Range("D349").Value = "VÄSTMANLANDS LÄN"
$ws.Range("E349").Value = "SKINNSKATTEBERG"
$ws.Range("F349").Value = "Sveaskog"
$ws.Range("G349").Value = 2
$ws.Range("H349").Value = 0
$ws.Range("I349").Value = 0
$ws.Range("J349").Value = 0
$ws.Range("K349").Value = 0
$ws.Range("L349").Value = 0
$ws.Range("M349").Value = 0
$ws.Range("N349").Value = 0
$ws.Range("O349").Value = 0
$ws.Range("P349").Value = 0
$ws.Range("Q349").Value = 0
$ws.Range("R349").WrapText = $true
